$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column C from Total_Profit to Average_Profit
$ws.Cells.Item(1, 3).Value = "Average_Profit"

# Find the last used row in column A (product_combinations)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $occurrences = $ws.Cells.Item($r, 2).Value()
    $totalProfit = $ws.Cells.Item($r, 3).Value()
    if ($occurrences -ne $null -and $occurrences -ne 0) {
        $ws.Cells.Item($r, 3).Value = $totalProfit / $occurrences
    }
}
